$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Clear cells that no longer hold data in their original positions ---
$ws.Range("C2").ClearContents()
$ws.Range("C9").ClearContents()
$ws.Range("A10").ClearContents()
$ws.Range("A11").ClearContents()
$ws.Range("C20").ClearContents()

# --- Row 2: Enter in filename ---
$ws.Range("A2").Value = "Enter in filename"
$ws.Range("B2").Value = "inputRacer.txt"

# --- Row 3: Display Items (new row, shifts everything below down) ---
$ws.Range("A3").Value = "Display Items"
$ws.Range("B3").Value = "P or D"
$ws.Range("C3").Value = "P is Tree, D is hash"

# --- Row 4: Display Menu (previously row 2) ---
$ws.Range("A4").Value = "Display Menu"
$ws.Range("B4").Value = "M"
$ws.Range("C4").Value = "Displays Menu"

# --- Row 5: Search by Name ---
$ws.Range("A5").Value = "Search by Name (BST Search)"
$ws.Range("B5").Value = 'S + "Pooh Bear"'
$ws.Range("C5").Value = "Shows Racer Result"

# --- Row 6: Search by ID ---
$ws.Range("A6").Value = "Search by ID (Hash Table Search)"
$ws.Range("B6").Value = 'F + "K0098983"'
$ws.Range("C6").Value = "Shows Racer Result"

# --- Row 7: Add a Racer ---
$ws.Range("A7").Value = "Add a Racer"
$ws.Range("B7").Value = "Name: Barack Obama "

# --- Row 8: ID ---
$ws.Range("B8").Value = "ID: A1234"

# --- Row 9: Circuit ---
$ws.Range("B9").Value = "Circuit: White House"

# --- Row 10: Date ---
$ws.Range("B10").Value = "Date: 01/21/13"

# --- Row 11: Laptime / Adds Barack Obama ---
$ws.Range("B11").Value = "Laptime: 100"
$ws.Range("C11").Value = "Adds Barack Obama"

# --- Row 12: Display Hash Table ---
$ws.Range("A12").Value = "Display Hash Table"
$ws.Range("B12").Value = "D"

# --- Row 13: Print in alphabetical ---
$ws.Range("A13").Value = "Print in alphabetical "
$ws.Range("B13").Value = "P"

# --- Row 14: Delete Barack Obama (replaces Remove by License) ---
$ws.Range("A14").Value = "Delete Barack Obama "
$ws.Range("B14").Value = 'R + "A1234"'
$ws.Range("C14").Value = "Removes Barack"

# --- Row 15: Display and Print to Show Gone ---
$ws.Range("A15").Value = "Display and Print to Show Gone"
$ws.Range("B15").Value = "D + P"

# --- Row 16: Undo delete ---
$ws.Range("A16").Value = "Undo delete"
$ws.Range("B16").Value = "Z"
$ws.Range("C16").Value = "barack back"

# --- Row 17: Display and Print to show its back ---
$ws.Range("A17").Value = "Display and Print to show its back"
$ws.Range("B17").Value = "D + P"

# --- Row 18: Show undo Z again ---
$ws.Range("A18").Value = "Show undo Z again"
$ws.Range("B18").Value = "Z"
$ws.Range("C18").Value = "Barack is gone again!"

# --- Row 19: Z again ---
$ws.Range("A19").Value = "Z again"
$ws.Range("B19").Value = "Z"
$ws.Range("C19").Value = '"No previous commands"'

# --- Row 20: Show statistics (previously row 12) ---
$ws.Range("A20").Value = "Show statistics"
$ws.Range("B20").Value = "T"

# --- Row 21: Show Indented Tree (previously row 13) ---
$ws.Range("A21").Value = "Show Indented Tree"
$ws.Range("B21").Value = "I"

# --- Row 22: Add Someone ---
$ws.Range("A22").Value = "Add Someone"
$ws.Range("B22").Value = "George Bush"

# --- Row 23: Texas ---
$ws.Range("B23").Value = "Texas"

# --- Row 24: A000 ---
$ws.Range("B24").Value = "A000"

# --- Row 25: Date: 01/21/13 Laptime: 100 ---
$ws.Range("B25").Value = "Date: 01/21/13 Laptime: 100"

# --- Row 26: Save hash table output to file (previously row 18) ---
$ws.Range("A26").Value = "Save hash table output to file"
$ws.Range("B26").Value = "O"

# --- Row 28: Q, Q, Q (previously row 20; row 27 left blank) ---
$ws.Range("A28").Value = "Q"
$ws.Range("B28").Value = "Q"
$ws.Range("C28").Value = "Q"

# --- Update selection to match the saved view state ---
$ws.Range("C16").Select() | Out-Null
